# Update cryptos price list (Price and Volume(1h) columns) to reflect
# the latest GitHub Actions scrape, including two coins that swapped
# ranking position (rows 38/39 and 45/46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.Value = '''29.081.33'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.13%  '

# Row 3
$c = $ws.Range("D3")
$c.Value = '''1.835.90'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.39%  '

# Row 4
$c = $ws.Range("D4")
$c.Value = '''1.0000'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$c = $ws.Range("D5")
$c.Value = '''244.64'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '

# Row 6
$c = $ws.Range("D6")
$c.Value = '''0.6353'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.08%  '

# Row 7
$c = $ws.Range("D7")
$c.Value = '''1.002'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.14%  '

# Row 8
$c = $ws.Range("D8")
$c.Value = '''0.07568'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +2.43%  '

# Row 9
$c = $ws.Range("D9")
$c.Value = '''0.2950'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.00%  '

# Row 10
$c = $ws.Range("D10")
$c.Value = '''22.93'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.10%  '

# Row 11
$c = $ws.Range("D11")
$c.Value = '''0.07753'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.99%  '

# Row 12
$c = $ws.Range("D12")
$c.Value = '''1.842.74'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.73%  '

# Row 13
$c = $ws.Range("D13")
$c.Value = '''5.008'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.01%  '

# Row 14
$c = $ws.Range("D14")
$c.Value = '''0.6715'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.34%  '

# Row 15
$c = $ws.Range("D15")
$c.Value = '''83.31'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.46%  '

# Row 16
$c = $ws.Range("D16")
$c.Value = '''0.000009601'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +5.71%  '

# Row 17
$c = $ws.Range("D17")
$c.Value = '''6.122'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.86%  '

# Row 18
$c = $ws.Range("D18")
$c.Value = '''29.121.85'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.26%  '

# Row 19
$c = $ws.Range("D19")
$c.Value = '''12.59'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.98%  '

# Row 20
$c = $ws.Range("D20")
$c.Value = '''226.81'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.72%  '

# Row 21
$ws.Range("E21").Value = '  +0.06%  '

# Row 22
$c = $ws.Range("D22")
$c.Value = '''7.221'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.77%  '

# Row 23
$c = $ws.Range("D23")
$c.Value = '''1.002'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '

# Row 24
$c = $ws.Range("D24")
$c.Value = '''160.78'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.93%  '

# Row 25
$c = $ws.Range("D25")
$c.Value = '''0.1404'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.62%  '

# Row 26
$c = $ws.Range("D26")
$c.Value = '''8.548'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.59%  '

# Row 27
$c = $ws.Range("D27")
$c.Value = '''17.98'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.02%  '

# Row 28
$ws.Range("E28").Value = '  +0.33%  '

# Row 29
$c = $ws.Range("D29")
$c.Value = '''4.128'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.83%  '

# Row 30
$c = $ws.Range("D30")
$c.Value = '''4.073'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.15%  '

# Row 31
$ws.Range("E31").Value = '  +0.29%  '

# Row 32
$c = $ws.Range("D32")
$c.Value = '''0.05411'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +3.29%  '

# Row 33
$c = $ws.Range("D33")
$c.Value = '''1.865'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.53%  '

# Row 34
$c = $ws.Range("D34")
$c.Value = '''0.7480'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.79%  '

# Row 35
$c = $ws.Range("D35")
$c.Value = '''1.143'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.64%  '

# Row 36
$c = $ws.Range("D36")
$c.Value = '''2.663'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.61%  '

# Row 37
$c = $ws.Range("D37")
$c.Value = '''1.236.63'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.41%  '

# Row 38
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D38")
$c.Value = '''2.765'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.57%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D39")
$c.Value = '''0.01795'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.76%  '

# Row 40
$c = $ws.Range("D40")
$c.Value = '''6.621'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +4.70%  '

# Row 41
$c = $ws.Range("D41")
$c.Value = '''0.9073'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.12%  '

# Row 42
$c = $ws.Range("D42")
$c.Value = '''1.002'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.15%  '

# Row 43
$c = $ws.Range("D43")
$c.Value = '''102.22'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.65%  '

# Row 44
$c = $ws.Range("D44")
$c.Value = '''1.985.13'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.38%  '

# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range("D45")
$c.Value = '''0.00000000124'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +3.37%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D46")
$c.Value = '''65.09'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.25%  '

# Row 47
$c = $ws.Range("D47")
$c.Value = '''0.5112'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.07%  '

# Row 48
$c = $ws.Range("D48")
$c.Value = '''0.4090'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +3.21%  '

# Row 49
$c = $ws.Range("D49")
$c.Value = '''9.113'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.33%  '

# Row 50
$ws.Range("E50").Value = '  +1.71%  '

# Row 51
$c = $ws.Range("D51")
$c.Value = '''0.05779'
$c.Style = "Normal"
